# Update column C ("Förändrad") date values from 2023-09-08 (45177) to
# 2023-09-09 (45178) for every data row in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45177) {
        $cell.Value2 = 45178
    }
}
